# Update Sheets via scheduled runner: refresh market-price derived columns
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1553.5
$ws.Range("I70").Value = 1433.25
$ws.Range("J70").Value = 1633.6666
$ws.Range("K70").Value = 4299.75
$ws.Range("L70").Value = 4900.9998
$ws.Range("M70").Value = -4029.75
$ws.Range("N70").Value = -5440.9998
$ws.Range("H73").Value = 1553.5
$ws.Range("I73").Value = 1433.25
$ws.Range("J73").Value = 1633.6666
$ws.Range("K73").Value = 4299.75
$ws.Range("L73").Value = 4900.9998
$ws.Range("M73").Value = -3363.75
$ws.Range("N73").Value = -6772.9998
$ws.Range("H112").Value = 9072.656000000001
$ws.Range("J112").Value = 9870.585999999999
$ws.Range("L112").Value = 29611.758
$ws.Range("N112").Value = -31827.758
$ws.Range("H129").Value = 1339.1052
$ws.Range("I129").Value = 586.25
$ws.Range("J129").Value = 1886.6364
$ws.Range("K129").Value = 1758.75
$ws.Range("L129").Value = 5659.9092
$ws.Range("M129").Value = 3241.25
$ws.Range("N129").Value = -15659.9092
$ws.Range("H132").Value = 3472.0303
$ws.Range("I132").Value = 3268.0625
$ws.Range("K132").Value = 9804.1875
$ws.Range("M132").Value = -7274.1875
$ws.Range("H137").Value = 3310.3713
$ws.Range("I137").Value = 3061.1924
$ws.Range("J137").Value = 4030.2222
$ws.Range("K137").Value = 9183.5772
$ws.Range("L137").Value = 12090.6666
$ws.Range("M137").Value = -6633.5772
$ws.Range("N137").Value = -17190.6666
$ws.Range("H138").Value = 2319.318
$ws.Range("I138").Value = 2860.3572
$ws.Range("K138").Value = 8581.071599999999
$ws.Range("M138").Value = -3441.071599999999
$ws.Range("H141").Value = 6895.8184
$ws.Range("I141").Value = 3466.6667
$ws.Range("K141").Value = 10400.0001
$ws.Range("M141").Value = -5220.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2657.2964
$ws.Range("I61").Value = 1941.5
$ws.Range("J61").Value = 4088.889
$ws.Range("K61").Value = 1941.5
$ws.Range("L61").Value = 4088.889
$ws.Range("M61").Value = -1729.5
$ws.Range("N61").Value = -4512.889
$ws.Range("H74").Value = 1723.4117
$ws.Range("I74").Value = 1889.5
$ws.Range("J74").Value = 1324.8
$ws.Range("K74").Value = 1889.5
$ws.Range("L74").Value = 1324.8
$ws.Range("M74").Value = -1015.5
$ws.Range("N74").Value = -3072.8
$ws.Range("H77").Value = 1723.4117
$ws.Range("I77").Value = 1889.5
$ws.Range("J77").Value = 1324.8
$ws.Range("K77").Value = 9447.5
$ws.Range("L77").Value = 6624
$ws.Range("M77").Value = -5079.5
$ws.Range("N77").Value = -15360
$ws.Range("H136").Value = 2657.2964
$ws.Range("I136").Value = 1941.5
$ws.Range("J136").Value = 4088.889
$ws.Range("K136").Value = 5824.5
$ws.Range("L136").Value = 12266.667
$ws.Range("M136").Value = -3274.5
$ws.Range("N136").Value = -17366.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 18165.334
$ws.Range("I82").Value = 5023
$ws.Range("J82").Value = 28679.2
$ws.Range("K82").Value = 5023
$ws.Range("L82").Value = 28679.2
$ws.Range("M82").Value = -4640
$ws.Range("N82").Value = -29445.2
$ws.Range("H85").Value = 18165.334
$ws.Range("I85").Value = 5023
$ws.Range("J85").Value = 28679.2
$ws.Range("K85").Value = 5023
$ws.Range("L85").Value = 28679.2
$ws.Range("M85").Value = -3697
$ws.Range("N85").Value = -31331.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 21785.928
$ws.Range("I4").Value = 5000.25
$ws.Range("J4").Value = 28500.2
$ws.Range("K4").Value = 5000.25
$ws.Range("L4").Value = 28500.2
$ws.Range("M4").Value = -4888.25
$ws.Range("N4").Value = -28724.2
$ws.Range("H31").Value = 4595.4414
$ws.Range("I31").Value = 1080.4584
$ws.Range("J31").Value = 13031.4
$ws.Range("K31").Value = 1080.4584
$ws.Range("L31").Value = 13031.4
$ws.Range("M31").Value = -785.4584
$ws.Range("N31").Value = -13621.4
$ws.Range("H34").Value = 4595.4414
$ws.Range("I34").Value = 1080.4584
$ws.Range("J34").Value = 13031.4
$ws.Range("K34").Value = 1080.4584
$ws.Range("L34").Value = 13031.4
$ws.Range("M34").Value = -878.4584
$ws.Range("N34").Value = -13435.4
$ws.Range("H58").Value = 1502
$ws.Range("I58").Value = 725
$ws.Range("J58").Value = 2123.6
$ws.Range("K58").Value = 725
$ws.Range("L58").Value = 2123.6
$ws.Range("M58").Value = -522
$ws.Range("N58").Value = -2529.6
$ws.Range("H107").Value = 620.86664
$ws.Range("I107").Value = 583.61536
$ws.Range("J107").Value = 649.35297
$ws.Range("K107").Value = 583.61536
$ws.Range("L107").Value = 649.35297
$ws.Range("M107").Value = 1336.38464
$ws.Range("N107").Value = -4489.35297
$ws.Range("H132").Value = 5953828
$ws.Range("I132").Value = 957.13043
$ws.Range("J132").Value = 33337034
$ws.Range("K132").Value = 2871.39129
$ws.Range("L132").Value = 100011102
$ws.Range("M132").Value = -341.39129
$ws.Range("N132").Value = -100016162
$ws.Range("H136").Value = 1502
$ws.Range("I136").Value = 725
$ws.Range("J136").Value = 2123.6
$ws.Range("K136").Value = 2175
$ws.Range("L136").Value = 6370.799999999999
$ws.Range("M136").Value = 375
$ws.Range("N136").Value = -11470.8
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8753106
$ws.Range("I4").Value = 6670515
$ws.Range("J4").Value = 10002661
$ws.Range("K4").Value = 20011545
$ws.Range("L4").Value = 30007983
$ws.Range("M4").Value = -20011433
$ws.Range("N4").Value = -30008207
$ws.Range("H5").Value = 1762.8572
$ws.Range("I5").Value = 672
$ws.Range("K5").Value = 2016
$ws.Range("M5").Value = -1904
$ws.Range("H12").Value = 185.91667
$ws.Range("J12").Value = 224.88889
$ws.Range("L12").Value = 674.6666700000001
$ws.Range("N12").Value = -1020.66667
$ws.Range("H114").Value = 997
$ws.Range("I114").Value = 276
$ws.Range("J114").Value = 1429.6
$ws.Range("K114").Value = 828
$ws.Range("L114").Value = 4288.799999999999
$ws.Range("M114").Value = 2426
$ws.Range("N114").Value = -10796.8
$ws.Range("H115").Value = 4177
$ws.Range("I115").Value = 3297
$ws.Range("J115").Value = 4910.3335
$ws.Range("K115").Value = 9891
$ws.Range("L115").Value = 14731.0005
$ws.Range("M115").Value = -8716
$ws.Range("N115").Value = -17081.0005
$ws.Range("H121").Value = 1142.725
$ws.Range("I121").Value = 273
$ws.Range("J121").Value = 1213.2433
$ws.Range("K121").Value = 819
$ws.Range("L121").Value = 3639.7299
$ws.Range("M121").Value = 491
$ws.Range("N121").Value = -6259.7299
$ws.Range("H129").Value = 1864.25
$ws.Range("I129").Value = 517.25
$ws.Range("J129").Value = 2133.65
$ws.Range("K129").Value = 1551.75
$ws.Range("L129").Value = 6400.950000000001
$ws.Range("M129").Value = 3448.25
$ws.Range("N129").Value = -16400.95
$ws.Range("H135").Value = 1762.8572
$ws.Range("I135").Value = 672
$ws.Range("K135").Value = 6048
$ws.Range("M135").Value = -3513

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 90009
$ws.Range("J25").Value = 90009
$ws.Range("L25").Value = 90009
$ws.Range("N25").Value = -91067
$ws.Range("H122").Value = 4086.4146
$ws.Range("J122").Value = 5774
$ws.Range("L122").Value = 17322
$ws.Range("N122").Value = -22222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 18600.4
$ws.Range("I2").Value = 5000
$ws.Range("J2").Value = 32200.8
$ws.Range("K2").Value = 5000
$ws.Range("L2").Value = 32200.8
$ws.Range("M2").Value = -4888
$ws.Range("N2").Value = -32424.8
$ws.Range("H136").Value = 2610.375
$ws.Range("I136").Value = 3354.889
$ws.Range("J136").Value = 1653.1428
$ws.Range("K136").Value = 10064.667
$ws.Range("L136").Value = 4959.428400000001
$ws.Range("M136").Value = -7514.667000000001
$ws.Range("N136").Value = -10059.4284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2320.913
$ws.Range("I136").Value = 2018.4642
$ws.Range("J136").Value = 2791.389
$ws.Range("K136").Value = 6055.392599999999
$ws.Range("L136").Value = 8374.167000000001
$ws.Range("M136").Value = -3505.392599999999
$ws.Range("N136").Value = -13474.167

Write-Host "Applied market-data refresh across sheets"